$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1
$wdFindAsk = 0

function Get-ParagraphIndexByHeading($doc, $headingText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $headingText) {
            return $i
        }
    }
    return -1
}

# 1) Update activation date
$d.Content.Find.Execute("Ativação: 01/01/2012", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Ativação: 01/01/2024", $wdReplaceAll)

# 2) Remove EQD from the course list
$d.Content.Find.Execute("Curso (semestre ideal): EA (6), EB (5), EQD (6), EQN (6)", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Curso (semestre ideal): EA (6), EB (5), EQN (6)", $wdReplaceAll)

# 3) "Programa" paragraph: collapse the 7 numbered items (currently separated by manual
#    line breaks) into a single run of text, by stripping the manual line break
#    characters (char 11) from that paragraph only.
$progHeadingIdx = Get-ParagraphIndexByHeading $d "Programa"
$progPara = $d.Paragraphs($progHeadingIdx + 1)
$progRng = $progPara.Range.Duplicate()
$progRng.Find.Execute([char]11, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceAll)

# 4) "Método:" description text
$d.Content.Find.Execute("A avaliação da disciplina será feita por meio de avaliações escritas individuais (provas) e avaliações de atividades em grupo (relatórios das aulas práticas e/ou trabalhos escritos e/ou apresentações de seminários).", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "A avaliação será feita por meio de duas provas (P1 e P2). A critério do professor, a avaliação poderá ser complementada por meio de trabalhos e/ou relatórios, valendo até 30% da nota das provas.", $wdReplaceAll)

# 5) "Critério:" description text
$d.Content.Find.Execute("A Média Final (MF) será calculada pela média entre todas as avaliações realizadas durante o semestre, sendo o conjunto das avaliações individuais correspondentes a 75% da composição de MF e o conjunto das avaliações em grupo correspondentes a 25% da composição de MF. Será aprovado o aluno que obtiver MF maior ou igual a cinco e frequência mínima de 70% no semestre.", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "A nota final (NF) será calculada pela média aritmética das provas. NF=(P1 +P2)/2.", $wdReplaceAll)

# 6) "Norma de recuperação:" - merge its two sentences (separated by a manual break) into
#    a single run, then update the wording. First remove the manual break that separates
#    the two old sentences, restricting the search to the range that starts right at
#    "Norma de recuperação:" and runs through the end of that paragraph (so the other
#    manual breaks used after "Método:" / "Critério:" remain untouched).
$avalHeadingIdx = Get-ParagraphIndexByHeading $d "Avaliação"
$avalPara = $d.Paragraphs($avalHeadingIdx + 1)
$normaFind = $avalPara.Range.Duplicate()
$normaFind.Find.Execute("Norma de recuperação:", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdFindAsk)
$normaRng = $d.Range($normaFind.Start, $avalPara.Range.End)
$normaRng.Find.Execute([char]11, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceAll)

$d.Content.Find.Execute("No período de Recuperação haverá horário previamente definido para resolução de dúvidas e será realizada uma avaliação escrita individual (Prova da Recuperação = PR), com conteúdo de todos os tópicos apresentados na disciplina durante o semestre.A Nota de Recuperação (NR) será dada pela média aritmética entre a Média do Semestre (MF) e a Prova da Recuperação (PR), sendo considerado aprovado o aluno que obtiver NR maior ou igual a cinco.", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.", $wdReplaceAll)

# 7) "Bibliografia" paragraph: rewrite as a single consolidated run of text (several
#    references are removed/renumbered and two entries are updated).
$biblioHeadingIdx = Get-ParagraphIndexByHeading $d "Bibliografia"
$biblioPara = $d.Paragraphs($biblioHeadingIdx + 1)
$biblioRng = $d.Range($biblioPara.Range.Start, $biblioPara.Range.End - 1)
$biblioRng.Text = "1) Skoog, D.A.; Holler, F.J. ; Nieman, T.A. Princípios de análise instrumental. 6a. ed. Porto Alegre: Bookman, 2009.2) KRUG, F.J. (org.) Métodos de preparo de amostras: fundamentos sobre métodos de preparo de amostras orgânicas e inorgânicas para análise elementar. 1. ed. Piracicaba: Edição do autor, 2008.3) COLLINS, C.H.; BRAGA, G.L.; BONATO, P.S. (Org.) Fundamentos de cromatografia. 1. ed. Campinas: Editora da UNICAMP, 2006.Bibliografia complementar1) CHRISTIAN, G.D. Analytical chemistry. 4. ed. Nova York: John Wiley & Sons, 1986.3) SILVERSTEIN, R.M.; WEBSTER, F.X.; KIEMLE, D.J. Identificação espectrométrica de compostos orgânicos. 7. ed. Rio de Janeiro: Livros Técnicos e Científicos, 2007.4) WILLARD, H.H.; MERRITE, L.; DEAB, J. Instrumentação analítica. Lisboa: Fundação Calouste Gulbekian, 1989."
